# feat: implement sell data creation in CreateSellView
# Append two newly-created sell records to the SellData sheet and
# correct the timestamp recorded for the previous (3rd) sale.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the timestamp of the existing row 4 (3rd sale) to match the
# value actually persisted by the (now cleaned-up) save flow.
$ws.Range("E4").Value = 45814.02324728009

# Row 5 - 4th sale created via CreateSellView
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "2AYB-3AYB-5AYB"
$ws.Range("C5").Value = "1-1-1"
$ws.Range("D5").Value = 27500
$ws.Range("E5").Value = 45815.77188895833
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 6 - 5th sale created via CreateSellView
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "1AYB-1AP-1M"
$ws.Range("C6").Value = "2-3-2"
$ws.Range("D6").Value = 45700
$ws.Range("E6").Value = 45815.77444396124
$ws.Range("E6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
